# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the panel_query_time ("time_taken") timestamps on the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:34:49.455892"
$dataSheet.Range("F3").Value = "2021-10-05 14:34:49.455899"
$dataSheet.Range("F4").Value = "2021-10-05 14:34:49.455903"
$dataSheet.Range("F5").Value = "2021-10-05 14:34:49.455905"
$dataSheet.Range("F6").Value = "2021-10-05 14:34:49.455908"
$dataSheet.Range("F7").Value = "2021-10-05 14:34:49.455911"
$dataSheet.Range("F8").Value = "2021-10-05 14:34:49.455913"
$dataSheet.Range("F9").Value = "2021-10-05 14:34:49.455916"
$dataSheet.Range("F10").Value = "2021-10-05 14:34:49.455919"
$dataSheet.Range("F11").Value = "2021-10-05 14:34:49.455921"
$dataSheet.Range("F12").Value = "2021-10-05 14:34:49.455924"
$dataSheet.Range("F13").Value = "2021-10-05 14:34:49.455927"
$dataSheet.Range("F14").Value = "2021-10-05 14:34:49.455929"
$dataSheet.Range("F15").Value = "2021-10-05 14:34:49.455932"
$dataSheet.Range("F16").Value = "2021-10-05 14:34:49.455935"
$dataSheet.Range("F17").Value = "2021-10-05 14:34:49.455937"
$dataSheet.Range("F18").Value = "2021-10-05 14:34:49.455940"
$dataSheet.Range("F19").Value = "2021-10-05 14:34:49.455943"
$dataSheet.Range("F20").Value = "2021-10-05 14:34:49.455945"
$dataSheet.Range("F21").Value = "2021-10-05 14:34:49.455948"
$dataSheet.Range("F22").Value = "2021-10-05 14:34:49.455950"
$dataSheet.Range("F23").Value = "2021-10-05 14:34:49.455953"
$dataSheet.Range("F24").Value = "2021-10-05 14:34:49.455956"
$dataSheet.Range("F25").Value = "2021-10-05 14:34:49.455959"

# --- Add the new "metadata" sheet, placed after "data" ---
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Metabolic renal disease"
$meta.Range("C2").Value = 198
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.26"
$meta.Range("E2").Value = "2021-09-09T06:53:39.643272Z"
$meta.Range("F2").Value = "2021-10-05 14:34:49.452511"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/198/?format=json"

# Match the header/index-cell formatting (bold, centered, bordered) used on the "data" sheet
$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
